$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 335.875
$ws.Cells.Item(96, 9).Value = 282.8
$ws.Cells.Item(96, 10).Value = 424.33334
$ws.Cells.Item(96, 11).Value = 848.4000000000001
$ws.Cells.Item(96, 12).Value = 1273.00002
$ws.Cells.Item(96, 13).Value = 524.5999999999999
$ws.Cells.Item(96, 14).Value = -4019.00002

$ws.Cells.Item(103, 8).Value = 9091561
$ws.Cells.Item(103, 9).Value = 482.63635
$ws.Cells.Item(103, 11).Value = 1447.90905
$ws.Cells.Item(103, 13).Value = -861.90905

$ws.Cells.Item(107, 8).Value = 2050.6667
$ws.Cells.Item(107, 9).Value = 2580
$ws.Cells.Item(107, 10).Value = 1389
$ws.Cells.Item(107, 11).Value = 2580
$ws.Cells.Item(107, 12).Value = 1389
$ws.Cells.Item(107, 13).Value = -660
$ws.Cells.Item(107, 14).Value = -5229

$ws.Cells.Item(116, 8).Value = 5044.8184
$ws.Cells.Item(116, 9).Value = 3695.5557
$ws.Cells.Item(116, 10).Value = 5978.923
$ws.Cells.Item(116, 11).Value = 3695.5557
$ws.Cells.Item(116, 12).Value = 5978.923
$ws.Cells.Item(116, 13).Value = -253.5556999999999
$ws.Cells.Item(116, 14).Value = -12862.923

$ws.Cells.Item(137, 8).Value = 1878.6305
$ws.Cells.Item(137, 9).Value = 1111.9143
$ws.Cells.Item(137, 10).Value = 4318.1816
$ws.Cells.Item(137, 11).Value = 3335.7429
$ws.Cells.Item(137, 12).Value = 12954.5448
$ws.Cells.Item(137, 13).Value = -785.7428999999997
$ws.Cells.Item(137, 14).Value = -18054.5448

$ws.Cells.Item(141, 8).Value = 2008
$ws.Cells.Item(141, 9).Value = 1854.5
$ws.Cells.Item(141, 10).Value = 2315
$ws.Cells.Item(141, 11).Value = 5563.5
$ws.Cells.Item(141, 12).Value = 6945
$ws.Cells.Item(141, 13).Value = -383.5
$ws.Cells.Item(141, 14).Value = -17305

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1834.0264
$ws.Cells.Item(2, 9).Value = 1220.6786
$ws.Cells.Item(2, 10).Value = 3551.4
$ws.Cells.Item(2, 11).Value = 1220.6786
$ws.Cells.Item(2, 12).Value = 3551.4
$ws.Cells.Item(2, 13).Value = -1107.6786
$ws.Cells.Item(2, 14).Value = -3777.4

$ws.Cells.Item(4, 8).Value = 287.5
$ws.Cells.Item(4, 9).Value = 275
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 275
$ws.Cells.Item(4, 12).Value = 300
$ws.Cells.Item(4, 13).Value = -159
$ws.Cells.Item(4, 14).Value = -532

$ws.Cells.Item(74, 8).Value = 1912.9706
$ws.Cells.Item(74, 9).Value = 858.1070999999999
$ws.Cells.Item(74, 10).Value = 6835.6665
$ws.Cells.Item(74, 11).Value = 858.1070999999999
$ws.Cells.Item(74, 12).Value = 6835.6665
$ws.Cells.Item(74, 13).Value = 15.89290000000005
$ws.Cells.Item(74, 14).Value = -8583.666499999999

$ws.Cells.Item(77, 8).Value = 1912.9706
$ws.Cells.Item(77, 9).Value = 858.1070999999999
$ws.Cells.Item(77, 10).Value = 6835.6665
$ws.Cells.Item(77, 11).Value = 4290.5355
$ws.Cells.Item(77, 12).Value = 34178.3325
$ws.Cells.Item(77, 13).Value = 77.46450000000004
$ws.Cells.Item(77, 14).Value = -42914.3325

$ws.Cells.Item(116, 8).Value = 1834.0264
$ws.Cells.Item(116, 9).Value = 1220.6786
$ws.Cells.Item(116, 10).Value = 3551.4
$ws.Cells.Item(116, 11).Value = 1220.6786
$ws.Cells.Item(116, 12).Value = 3551.4
$ws.Cells.Item(116, 13).Value = 1073.3214
$ws.Cells.Item(116, 14).Value = -8139.4

$ws.Cells.Item(135, 8).Value = 37679.332
$ws.Cells.Item(135, 10).Value = 37679.332
$ws.Cells.Item(135, 12).Value = 37679.332
$ws.Cells.Item(135, 14).Value = -47819.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1834.0264
$ws.Cells.Item(3, 9).Value = 1220.6786
$ws.Cells.Item(3, 10).Value = 3551.4
$ws.Cells.Item(3, 11).Value = 1220.6786
$ws.Cells.Item(3, 12).Value = 3551.4
$ws.Cells.Item(3, 13).Value = -1106.6786
$ws.Cells.Item(3, 14).Value = -3779.4

$ws.Cells.Item(134, 8).Value = 2099.4211
$ws.Cells.Item(134, 9).Value = 1828.6061
$ws.Cells.Item(134, 10).Value = 3886.8
$ws.Cells.Item(134, 11).Value = 5485.8183
$ws.Cells.Item(134, 12).Value = 11660.4
$ws.Cells.Item(134, 13).Value = -2950.8183
$ws.Cells.Item(134, 14).Value = -16730.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 32.866665
$ws.Cells.Item(7, 9).Value = 30.916666
$ws.Cells.Item(7, 10).Value = 40.666668
$ws.Cells.Item(7, 11).Value = 30.916666
$ws.Cells.Item(7, 12).Value = 40.666668
$ws.Cells.Item(7, 13).Value = 82.08333400000001
$ws.Cells.Item(7, 14).Value = -266.666668

$ws.Cells.Item(45, 8).Value = 8500
$ws.Cells.Item(45, 10).Value = 8500
$ws.Cells.Item(45, 12).Value = 8500
$ws.Cells.Item(45, 14).Value = -9686

$ws.Cells.Item(51, 8).Value = 9124.833000000001
$ws.Cells.Item(51, 10).Value = 9124.833000000001
$ws.Cells.Item(51, 12).Value = 9124.833000000001
$ws.Cells.Item(51, 14).Value = -10596.833

$ws.Cells.Item(61, 8).Value = 9124.833000000001
$ws.Cells.Item(61, 10).Value = 9124.833000000001
$ws.Cells.Item(61, 12).Value = 9124.833000000001
$ws.Cells.Item(61, 14).Value = -9820.833000000001

$ws.Cells.Item(132, 8).Value = 1822.5769
$ws.Cells.Item(132, 9).Value = 2388.875
$ws.Cells.Item(132, 10).Value = 916.5
$ws.Cells.Item(132, 11).Value = 7166.625
$ws.Cells.Item(132, 12).Value = 2749.5
$ws.Cells.Item(132, 13).Value = -4636.625
$ws.Cells.Item(132, 14).Value = -7809.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 4699.8887
$ws.Cells.Item(51, 9).Value = 1199.75
$ws.Cells.Item(51, 10).Value = 7500
$ws.Cells.Item(51, 11).Value = 3599.25
$ws.Cells.Item(51, 12).Value = 22500
$ws.Cells.Item(51, 13).Value = -3139.25
$ws.Cells.Item(51, 14).Value = -23420

$ws.Cells.Item(103, 8).Value = 860

$ws.Cells.Item(131, 8).Value = 950.65
$ws.Cells.Item(131, 10).Value = 963.62244
$ws.Cells.Item(131, 12).Value = 2890.86732
$ws.Cells.Item(131, 14).Value = -12970.86732

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 7972.6
$ws.Cells.Item(41, 9).Value = 6000
$ws.Cells.Item(41, 10).Value = 8465.75
$ws.Cells.Item(41, 11).Value = 6000
$ws.Cells.Item(41, 12).Value = 8465.75
$ws.Cells.Item(41, 13).Value = -5645
$ws.Cells.Item(41, 14).Value = -9175.75

$ws.Cells.Item(113, 8).Value = 1582.5714
$ws.Cells.Item(113, 9).Value = 1553.1875
$ws.Cells.Item(113, 10).Value = 1676.6
$ws.Cells.Item(113, 11).Value = 1553.1875
$ws.Cells.Item(113, 12).Value = 1676.6
$ws.Cells.Item(113, 13).Value = 616.8125
$ws.Cells.Item(113, 14).Value = -6016.6

$ws.Cells.Item(132, 8).Value = 2547.36
$ws.Cells.Item(132, 9).Value = 2328.2104
$ws.Cells.Item(132, 10).Value = 3241.3333
$ws.Cells.Item(132, 11).Value = 6984.6312
$ws.Cells.Item(132, 12).Value = 9723.999899999999
$ws.Cells.Item(132, 13).Value = -4454.6312
$ws.Cells.Item(132, 14).Value = -14783.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(57, 8).Value = 4520.5
$ws.Cells.Item(57, 10).Value = 5000
$ws.Cells.Item(57, 12).Value = 5000
$ws.Cells.Item(57, 14).Value = -6132

$ws.Cells.Item(58, 8).Value = 1593
$ws.Cells.Item(58, 9).Value = 1593
$ws.Cells.Item(58, 11).Value = 1593
$ws.Cells.Item(58, 13).Value = -1333

$ws.Cells.Item(132, 8).Value = 4603.483
$ws.Cells.Item(132, 9).Value = 4532.4546
$ws.Cells.Item(132, 10).Value = 4826.7144
$ws.Cells.Item(132, 11).Value = 13597.3638
$ws.Cells.Item(132, 12).Value = 14480.1432
$ws.Cells.Item(132, 13).Value = -11067.3638
$ws.Cells.Item(132, 14).Value = -19540.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2982.5557
$ws.Cells.Item(107, 9).Value = 1680.2727
$ws.Cells.Item(107, 10).Value = 5029
$ws.Cells.Item(107, 11).Value = 5040.8181
$ws.Cells.Item(107, 12).Value = 15087
$ws.Cells.Item(107, 13).Value = -3120.8181
$ws.Cells.Item(107, 14).Value = -18927

$ws.Cells.Item(113, 8).Value = 46447.91
$ws.Cells.Item(113, 9).Value = 59150.06
$ws.Cells.Item(113, 10).Value = 3260.6
$ws.Cells.Item(113, 11).Value = 177450.18
$ws.Cells.Item(113, 12).Value = 9781.799999999999
$ws.Cells.Item(113, 13).Value = -175280.18
$ws.Cells.Item(113, 14).Value = -14121.8

$ws.Cells.Item(132, 8).Value = 4645.375
$ws.Cells.Item(132, 9).Value = 10701.6
$ws.Cells.Item(132, 10).Value = 1892.5454
$ws.Cells.Item(132, 11).Value = 32104.8
$ws.Cells.Item(132, 12).Value = 5677.6362
$ws.Cells.Item(132, 13).Value = -29574.8
$ws.Cells.Item(132, 14).Value = -10737.6362
